$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns (I0, IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (border/bold/centered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$iValues = @{
    2 = 9; 3 = 9; 4 = 9; 5 = 9; 6 = 4; 7 = 8; 8 = 7; 9 = 9; 10 = 11
    11 = 5; 12 = 7; 13 = 8; 14 = 6; 15 = 6; 16 = 6; 17 = 6; 18 = 8; 19 = 7
    20 = 6; 21 = 8; 22 = 9; 23 = 8; 24 = 4; 25 = 5; 26 = 8; 27 = 5; 28 = 8
    29 = 8; 30 = 5; 31 = 7; 32 = 7; 33 = 4; 34 = 6; 35 = 8; 36 = 9; 37 = 7
    38 = 9; 39 = 8; 40 = 6; 41 = 6; 42 = 5; 43 = 7; 44 = 7; 45 = 3; 46 = 7
    47 = 6; 48 = 6; 49 = 6
}

$jValues = @{
    2 = 9; 3 = 9; 4 = 9; 5 = 9; 6 = 4; 7 = 8; 8 = 8; 9 = 9; 10 = 11
    11 = 5; 12 = 7; 13 = 8; 14 = 6; 15 = 7; 16 = 6; 17 = 7; 18 = 8; 19 = 8
    20 = 6; 21 = 8; 22 = 9; 23 = 8; 24 = 4; 25 = 5; 26 = 8; 27 = 6; 28 = 8
    29 = 8; 30 = 5; 31 = 8; 32 = 7; 33 = 5; 34 = 6; 35 = 8; 36 = 9; 37 = 7
    38 = 9; 39 = 8; 40 = 6; 41 = 6; 42 = 5; 43 = 7; 44 = 7; 45 = 4; 46 = 7
    47 = 6; 48 = 6; 49 = 6
}

for ($row = 2; $row -le 49; $row++) {
    $ws.Cells.Item($row, 9).Value = $iValues[$row]
    $ws.Cells.Item($row, 10).Value = $jValues[$row]
}
